# Auto-generated Excel COM-interop edit script
# Applies updated market-price derived values to the Leve profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW) per scheduled runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1921.4595
$ws.Range("J17").Value = 1921.4595
$ws.Range("L17").Value = 5764.3785
$ws.Range("N17").Value = -6100.3785

$ws.Range("H62").Value = 5665
$ws.Range("I62").Value = 6128.1816
$ws.Range("J62").Value = 3966.6667
$ws.Range("K62").Value = 6128.1816
$ws.Range("L62").Value = 3966.6667
$ws.Range("M62").Value = -5504.1816
$ws.Range("N62").Value = -5214.6667

$ws.Range("H65").Value = 5665
$ws.Range("I65").Value = 6128.1816
$ws.Range("J65").Value = 3966.6667
$ws.Range("K65").Value = 30640.908
$ws.Range("L65").Value = 19833.3335
$ws.Range("M65").Value = -27520.908
$ws.Range("N65").Value = -26073.3335

$ws.Range("H129").Value = 1022.6022
$ws.Range("J129").Value = 969.4815
$ws.Range("L129").Value = 2908.4445
$ws.Range("N129").Value = -12908.4445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 77910.81
$ws.Range("I45").Value = 111950.82
$ws.Range("J45").Value = 3022.8
$ws.Range("K45").Value = 111950.82
$ws.Range("L45").Value = 3022.8
$ws.Range("M45").Value = -111573.82
$ws.Range("N45").Value = -3776.8

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H97").Value = 1084.6957
$ws.Range("I97").Value = 965.1053000000001
$ws.Range("K97").Value = 965.1053000000001
$ws.Range("M97").Value = -469.1053000000001

$ws.Range("H122").Value = 2305.5789
$ws.Range("I122").Value = 2089.8333
$ws.Range("J122").Value = 2675.4285
$ws.Range("K122").Value = 6269.499899999999
$ws.Range("L122").Value = 8026.2855
$ws.Range("M122").Value = -3819.499899999999
$ws.Range("N122").Value = -12926.2855

$ws.Range("H132").Value = 62503876
$ws.Range("I132").Value = 250001000
$ws.Range("J132").Value = 4832.3335
$ws.Range("K132").Value = 750003000
$ws.Range("L132").Value = 14497.0005
$ws.Range("M132").Value = -750000470
$ws.Range("N132").Value = -19557.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2023.6471
$ws.Range("I86").Value = 1861.6923
$ws.Range("J86").Value = 2550
$ws.Range("K86").Value = 1861.6923
$ws.Range("L86").Value = 2550
$ws.Range("M86").Value = -738.6922999999999
$ws.Range("N86").Value = -4796

$ws.Range("H89").Value = 2023.6471
$ws.Range("I89").Value = 1861.6923
$ws.Range("J89").Value = 2550
$ws.Range("K89").Value = 9308.461499999999
$ws.Range("L89").Value = 12750
$ws.Range("M89").Value = -3692.461499999999
$ws.Range("N89").Value = -23982

$ws.Range("H94").Value = 1392.4073
$ws.Range("I94").Value = 1399.9412
$ws.Range("J94").Value = 1379.6
$ws.Range("K94").Value = 1399.9412
$ws.Range("L94").Value = 1379.6
$ws.Range("M94").Value = -948.9412
$ws.Range("N94").Value = -2281.6

$ws.Range("H99").Value = 2093.303
$ws.Range("I99").Value = 1958.4814
$ws.Range("J99").Value = 2700
$ws.Range("K99").Value = 1958.4814
$ws.Range("L99").Value = 2700
$ws.Range("M99").Value = -460.4813999999999
$ws.Range("N99").Value = -5696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 522
$ws.Range("I7").Value = 182.75
$ws.Range("J7").Value = 861.25
$ws.Range("K7").Value = 182.75
$ws.Range("L7").Value = 861.25
$ws.Range("M7").Value = -69.75
$ws.Range("N7").Value = -1087.25

$ws.Range("H31").Value = 5213382
$ws.Range("I31").Value = 2487.158
$ws.Range("J31").Value = 7413537.5
$ws.Range("K31").Value = 2487.158
$ws.Range("L31").Value = 7413537.5
$ws.Range("M31").Value = -2192.158
$ws.Range("N31").Value = -7414127.5

$ws.Range("H34").Value = 5213382
$ws.Range("I34").Value = 2487.158
$ws.Range("J34").Value = 7413537.5
$ws.Range("K34").Value = 2487.158
$ws.Range("L34").Value = 7413537.5
$ws.Range("M34").Value = -2285.158
$ws.Range("N34").Value = -7413941.5

$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27372

$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -86856

$ws.Range("H132").Value = 54054.48
$ws.Range("I132").Value = 1520.8889
$ws.Range("J132").Value = 159121.67
$ws.Range("K132").Value = 4562.6667
$ws.Range("L132").Value = 477365.01
$ws.Range("M132").Value = -2032.6667
$ws.Range("N132").Value = -482425.01

$ws.Range("H141").Value = 3199
$ws.Range("J141").Value = 3199
$ws.Range("L141").Value = 3199
$ws.Range("N141").Value = -13559

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5719.7
$ws.Range("I5").Value = 6505.2354
$ws.Range("J5").Value = 1268.3334
$ws.Range("K5").Value = 19515.7062
$ws.Range("L5").Value = 3805.0002
$ws.Range("M5").Value = -19403.7062
$ws.Range("N5").Value = -4029.0002

$ws.Range("H97").Value = 29662.785
$ws.Range("J97").Value = 29662.785
$ws.Range("L97").Value = 88988.355
$ws.Range("N97").Value = -89980.355

$ws.Range("H113").Value = 13657
$ws.Range("J113").Value = 1309
$ws.Range("L113").Value = 3927
$ws.Range("N113").Value = -8267

$ws.Range("H131").Value = 852.38
$ws.Range("I131").Value = 626.6667
$ws.Range("J131").Value = 859.3608400000001
$ws.Range("K131").Value = 1880.0001
$ws.Range("L131").Value = 2578.08252
$ws.Range("M131").Value = 3159.9999
$ws.Range("N131").Value = -12658.08252

$ws.Range("H135").Value = 5719.7
$ws.Range("I135").Value = 6505.2354
$ws.Range("J135").Value = 1268.3334
$ws.Range("K135").Value = 58547.11859999999
$ws.Range("L135").Value = 11415.0006
$ws.Range("M135").Value = -56012.11859999999
$ws.Range("N135").Value = -16485.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 7925
$ws.Range("J21").Value = 7925
$ws.Range("L21").Value = 7925
$ws.Range("N21").Value = -8271

$ws.Range("H30").Value = 7925
$ws.Range("J30").Value = 7925
$ws.Range("L30").Value = 7925
$ws.Range("N30").Value = -8135

$ws.Range("H43").Value = 6999.2856
$ws.Range("I43").Value = 3499.25
$ws.Range("J43").Value = 11666
$ws.Range("K43").Value = 3499.25
$ws.Range("L43").Value = 11666
$ws.Range("M43").Value = -3348.25
$ws.Range("N43").Value = -11968

$ws.Range("H46").Value = 24082.6
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 25445.643
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 25445.643
$ws.Range("M46").Value = -4844
$ws.Range("N46").Value = -25757.643

$ws.Range("H122").Value = 2265.3333
$ws.Range("I122").Value = 1771.4286
$ws.Range("J122").Value = 3994
$ws.Range("K122").Value = 5314.2858
$ws.Range("L122").Value = 11982
$ws.Range("M122").Value = -2864.2858
$ws.Range("N122").Value = -16882

$ws.Range("H123").Value = 11930
$ws.Range("J123").Value = 11930
$ws.Range("L123").Value = 11930
$ws.Range("N123").Value = -16830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 30667.834
$ws.Range("J2").Value = 36501.4
$ws.Range("L2").Value = 36501.4
$ws.Range("N2").Value = -36725.4

$ws.Range("H22").Value = 971.4286
$ws.Range("I22").Value = 960
$ws.Range("K22").Value = 960
$ws.Range("M22").Value = -665

$ws.Range("H27").Value = 971.4286
$ws.Range("I27").Value = 960
$ws.Range("K27").Value = 960
$ws.Range("M27").Value = -853

$ws.Range("H122").Value = 93264
$ws.Range("I122").Value = 93264
$ws.Range("K122").Value = 279792
$ws.Range("M122").Value = -277342
